$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet "My Series" -> "Data"
$ws.Name = "Data"

# 2. Update the hidden CDM metadata blob stored in the cell comment on A1.
$commentText = "Jx0AAB+LCAAAAAAAAAOlWVtvG8cV/isLPrVAyV0ubUcWRhvwJoUoKQokVVl5KYa7I3Gq5S67MyuJbynQIkWaoigKp0iv6FOKAnWNNgFSu5f/Eliy+9S/0DOXvZFUzFUNwdo55zszZ86c24zQu9dz37gkEaNhsFep16yKQQI39GhwvleJ+Vm1/qjyroO61y7xj3CE54QD2ACpgO1eM7pXmXG+2DXNq6ur2lWjFkbnpm1ZdfPJoD92Z2SOqzRgHAcuqaRS3tulKg5qe/MB4djDHCvJvUpv3Ku1CXU7QBvgAJ+TqNaKGQ0IY92AU04JE5IRwZy0O4PvqI05du1RrY7MNXqGbMXU9xSugFR0jYNlyYTOiWNb9Z2qtVNtWJO6vWvv7NqN2o7deD8RTIGojxkfk+iSupIw5ni+kOLWTsOq2zZIIXMjCObKDOCgoe+NyCVlxGsT32elLGLqA2y6HHZdzpgWMnOyeqL7q3AQ4cVsQrlPyqkxGrSMeaB1ySZx0H4YERfsdy+VDsnVMNJmnSz6wJ3MaMSXHbwsPdcxI9FwIYxUTtRBnTDgTZ9E/HgBZ008cAVgODyKCTLvYGZCHcpc+KZBTDznDPssL1RgopMwumAL7JJDiGNTzHEV+CH2wOE4ZZy62aJrDHQUhQuYERZvhb63D7Nq8AZGOnMvABOLZVtheJFpt4mJ5KnK84UznWOewNfoaDwLr4aBvxzHU+ZGdEq8TitBb+QhEZBauh0zHs5Bi4yEFC1HWcI/CMBVMuoQl86xf+SDEZmIjSIBNWMenlHeDv14HrBEpxUqOoEdTch1usN0jIZwuIEwehj0ggSvzLyRVRQYhVfpmusMaYQcucnc5LjXGavgDtCS41vnyBMRu9ynPtSH/FnkqEWvGM8I4RtdQnGQSIX7ouI4reVhPJ9CeE0hxi7lqgyZGR+Bn4Kvg16OBVWkKn8mlrUrf0CPlI26gXc3LmEiWC63llMH3goJwZ78lo+DC6CeUD47bCZ72cBBygJ34td5CCJ34eOlJKdWytNQL3D92CMqIfSCM+miQjd1qHey0RqpDzHuIBwsJ8sF5GVGdzl87FWgUu8yHkEvUHHcMA54tBSZA5ka+jYZFk8DuQD2t5Y5i8j3Y2hBlvtx4LZDb/vVPGWd44Dy7TUM40ilw+1FpPVEZoxZh4gcI5P+1vJumT2xqBR8HpB5GFB3e2uDkYX23j02wpKo2lqCqPjaGu9DXVdlT8T61mIR9I9Q6Eot02QsdKl0Vh0eXk7evCNkOuQMxz70bhxK7Hmae1fJqMkuVjF5EjqO/CQDOqIzZtAau9685kLzINq/mhvOBcGEjvRkjMw8XnRALukG530cnMfQY6R5ZZWe5l9RHycRDpjYTtpSrKTizSCU5CnV6jgqeQ1j6QgqeYXAReYKDk3IfBFG2B+AYei+djvdL0ErMsB8pkdQ23ziJkY2M9FUqqhZovjbYLJIqW2IgNdpcoUoQWIvqgnPMBkNiV0OICz9NvbpNFJZNSnlm3hwYFlzmORfsbmSjWJyBnARg+r7bbIUHUg20HTpsvWEoRxYJFJnPHqwYz+0GjZ0NmKM5I5HBPtGF4KZE6MXXBLG5yC2a4wIox58UezvGu+RKaFQBKWJdBkqLZ2XQ/tJnpeqNEHfIqUIgH7jnEIZWQemnEzAOSU48pc5oNpqP3QBd/vjf9/85vmrF5/dfvT0zRc//O8/fvXqnz+7efYj+Lj9699uPv6l2qYCowme+kQqNGnt7FiNB+BnKQkJ45qyNfZil0va6ansiNMx0hc7OWh3e+2Dfkvmk5SYiKuSYoo74zKMs+FYbUIuJI/UTDxBQZxJkp/0uMDNlShH3OEuSRGd598lqGzx+uVnr1/++U5pbbCs16o/fvywWrff2orBbbi+hktbsX6hBgjwg6r1sGrbOfAKBo1UAUjt1POcRt16bNlw5U5zuZc68ibQKkvPNMHn5oqcIrVVe5S6QH6cMKXjTyBEUrYKhdxAu+gXP3nzl6cFlLauphRnAeVkGyMWM5OBnPpwNDHGw+NRu2tMumPhJxkvh1OTfw1Yr57GU8GpgiDG/rcMKO5QzIwK3IgqRnhmEOzOjCVEYi4OC862iaoWuueUq1oeRGG8UCeSE8ioG5BpNtkosSHXSJ6051rSyVgb4ErXm79/vklAb6STNbTp+0eehgocRcrxddR++q9XX3746sWL2+c/v/nyB4UZ9DrpswD4OURTfpi6PaQ8XW9WKOhkLI15YX03V180UVykjkIacObUH8k7lB4hEK2L2eRv1JtDyZMTS3sBfYWC3sOse811YDuHyCwSQM8FhmobZnfPlKByeGbX//z2d7e//vz2k+dvPvzTzUd/vPn4k9cvf//m2R9U1N0+fX7702c6y68WAqmLuNGqJtCQ7yOuIaLRELXb+OqDXxhByA1oOYxYZqSvPvg0N5lQVDYn2czQ0qWKFFVYg+aFhZyRUyXVoSCXiqgGoC1KWCNF6CIWLqibLfJ+VUwl4k4yvtGbVGNGjBC6qW/CTorgTHhbOS2iSurRO5ZdtzVXaSO2MMUsZ/oDP5xCk5Ew5APECqQg9fUCGVaud9Aftpr9DKKUGEYeiYQbqg+UtJSipPRYMkpcLUcBLjR+buyLN6M12DornTmXxkz9+HLW9ET62/xUUUCgdhxFqiEK9Fv+OF5AM5w80d3Nl8+Wuf73UPWq+Y44G/c6RT6Mc1wohEW2IEi+TE2apdJUj4l3HtXOHgrTZEPgFZ46wRz6uV51WpfQV0amyDvdKAqjjckn4ySwAXTSkFHMzOIpRp6p6rq97KwSQpLw0g9189M7DDvEJ7zcW7aZSQ/Cy3vLwtmXFe2xoe9pY5a7eqRmySbIP+gLR/l/3/OVszWjCBor8QBY+gE+ubiO4L5bUhu1FSkoboCwun4h36cR409EJtBfinKaUk5Vh/pEXLjUhxyfOo2HigAAMz+7WVAzCV2u/qwS+n06pyWvhVYS38VJwJaLhWrheuU8RZSWQ3INDWZuBkiK0+9B2VDvKGVmUw4LuTSVF2+XjJ7PeFnF3pli4pGpVXWnxK4+8Kyd6mNCGtV6Hf7Hrm1b1kPx8qknh8xByVXJRczkwLI/dzr/A5QeT1InHQAA"
$ws.Range("A1").Comment.Text($commentText) | Out-Null

# 3. Change the custom number format (numFmtId 166) used by B27:B36 from "0" to "###0"
$ws.Range("B27:B36").NumberFormat = "###0"

# 4. Update the label in A11
$ws.Range("A11").Value = "Function Information"

# 5. Update the Kurtosis value in B21
$ws.Range("B21").Value = 0.2499825759175085
